$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# uploader id record added: stash the uploader ids for the music entry in row 1
$ws.Range("B1").Value = "albida"
$ws.Range("C1").Value = "albida8"

# Comments recording attribution / uploader ids
$ws.Range("A1").AddComment("Text") | Out-Null
$ws.Range("B1").AddComment("1193015946") | Out-Null
$ws.Range("C1").AddComment("1193015946") | Out-Null
